$d = $word.ActiveDocument

function Replace-InRange($range, $find, $replace) {
    $range.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# Paragraph 1: "English / Portuguese / French / Thai / Vietnamese / Spanish"
$p = $d.Paragraphs.Item(1).Range
Replace-InRange $p "English" "Anglais"
$p = $d.Paragraphs.Item(1).Range
Replace-InRange $p " / Portuguese / French / Thai / Vietnamese / Spanish" " / portugais / français / thaïlandais / vietnamien / espagnol"

# Paragraph 3: "English" (standalone)
$p = $d.Paragraphs.Item(3).Range
Replace-InRange $p "English" "Anglais"

# Paragraph 5 (table cell): "Brief:"
$p = $d.Paragraphs.Item(5).Range
Replace-InRange $p "Brief" "Résumé"
$p = $d.Paragraphs.Item(5).Range
Replace-InRange $p ":" " :"

# Paragraph 6 (table cell): brief description
$p = $d.Paragraphs.Item(6).Range
Replace-InRange $p "It will be sent" "Il sera envoyé"

# Paragraph 8 (table cell): "Target audience:"
$p = $d.Paragraphs.Item(8).Range
Replace-InRange $p "Target audience" "Public cible"
$p = $d.Paragraphs.Item(8).Range
Replace-InRange $p ":" " :"

# Paragraph 14: "Uh oh! Your documents couldn't be verified"
$p = $d.Paragraphs.Item(14).Range
Replace-InRange $p "Uh oh! Your documents couldn’t be verified" "Oups ! Vos documents n'ont pas pu être vérifiés."

# Paragraph 16: "Hi [PARTNER NAME], "
$p = $d.Paragraphs.Item(16).Range
Replace-InRange $p "Hi " "Salut "
$p = $d.Paragraphs.Item(16).Range
Replace-InRange $p "[PARTNER NAME]" "[NOM DU PARTENAIRE]"

# Paragraph 17: "We regret to inform you..."
$p = $d.Paragraphs.Item(17).Range
Replace-InRange $p "We regret to inform you that your documents have failed our verification process as we found the following issues with them: " "Nous avons le regret de vous informer que vos documents n'ont pas été vérifiés pour les raisons suivantes : "

# Paragraph 18: "A copy of your vaccination certificate: Document is unclear"
$p = $d.Paragraphs.Item(18).Range
Replace-InRange $p "A copy of your vaccination certificate" "Copie de votre certificat de vaccination"
$p = $d.Paragraphs.Item(18).Range
Replace-InRange $p ": Document is unclear" ": document illisible"

# Paragraph 20: "Please resubmit the documents above by [DD Mmm YYYY] so we can proceed with the necessary arrangements."
$p = $d.Paragraphs.Item(20).Range
Replace-InRange $p "Please resubmit the documents above by " "Veuillez renvoyer les documents ci-dessus avant le "
$p = $d.Paragraphs.Item(20).Range
Replace-InRange $p "DD Mmm YYYY" "JJ Mmm AAAA"
$p = $d.Paragraphs.Item(20).Range
Replace-InRange $p " so we can proceed with the necessary arrangements." " afin que nous puissions prendre les dispositions nécessaires."

# Paragraph 21: "If you have any questions, please contact us via live chat or WhatsApp. "
$p = $d.Paragraphs.Item(21).Range
Replace-InRange $p "If you have any questions, please contact us via " "Si vous avez des questions, veuillez nous contacter par "
$p = $d.Paragraphs.Item(21).Range
Replace-InRange $p "live chat" "chat en direct"
$p = $d.Paragraphs.Item(21).Range
Replace-InRange $p " or " " ou sur "

# Paragraph 22: "If you have any questions, please contact your country manager, [NAME], at [EMAIL ADDRESS] or [WHATSAPP NO] (WhatsApp). "
$p = $d.Paragraphs.Item(22).Range
Replace-InRange $p "If you have any questions, please contact your country manager, " "Si vous avez des questions, veuillez contacter votre responsable national, "
$p = $d.Paragraphs.Item(22).Range
Replace-InRange $p ", at " ", à l'adresse "
$p = $d.Paragraphs.Item(22).Range
Replace-InRange $p " or " " ou au"

# Comment 0: "choose either one"
$c = $d.Comments.Item(1)
$c.Range.Text = "choisissez l'un ou l'autre"
